$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 246.82143
$ws.Range("I33").Value = 248.48148
$ws.Range("K33").Value = 248.48148
$ws.Range("M33").Value = -19.48148

$ws.Range("H58").Value = 1573.1578
$ws.Range("I58").Value = 251
$ws.Range("J58").Value = 3042.2222
$ws.Range("K58").Value = 753
$ws.Range("L58").Value = 9126.6666
$ws.Range("M58").Value = -603
$ws.Range("N58").Value = -9426.6666

$ws.Range("H98").Value = 9930.352999999999
$ws.Range("I98").Value = 7462.6924
$ws.Range("J98").Value = 17950.25
$ws.Range("K98").Value = 7462.6924
$ws.Range("L98").Value = 17950.25
$ws.Range("M98").Value = -5964.6924
$ws.Range("N98").Value = -20946.25

$ws.Range("H122").Value = 9930.352999999999
$ws.Range("I122").Value = 7462.6924
$ws.Range("J122").Value = 17950.25
$ws.Range("K122").Value = 22388.0772
$ws.Range("L122").Value = 53850.75
$ws.Range("M122").Value = -19938.0772
$ws.Range("N122").Value = -58750.75

$ws.Range("H129").Value = 952.5
$ws.Range("I129").Value = 518.6667
$ws.Range("J129").Value = 1029.0588
$ws.Range("K129").Value = 1556.0001
$ws.Range("L129").Value = 3087.1764
$ws.Range("M129").Value = 3443.9999
$ws.Range("N129").Value = -13087.1764

$ws.Range("H132").Value = 1659.174
$ws.Range("J132").Value = 4003.875
$ws.Range("L132").Value = 12011.625
$ws.Range("N132").Value = -17071.625

$ws.Range("H137").Value = 1507.4286
$ws.Range("I137").Value = 1348.6562
$ws.Range("J137").Value = 1806.2941
$ws.Range("K137").Value = 4045.9686
$ws.Range("L137").Value = 5418.8823
$ws.Range("M137").Value = -1495.9686
$ws.Range("N137").Value = -10518.8823

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5228.625
$ws.Range("I63").Value = 4771.5
$ws.Range("J63").Value = 6600
$ws.Range("K63").Value = 4771.5
$ws.Range("L63").Value = 6600
$ws.Range("M63").Value = -4085.5
$ws.Range("N63").Value = -7972

$ws.Range("H66").Value = 5228.625
$ws.Range("I66").Value = 4771.5
$ws.Range("J66").Value = 6600
$ws.Range("K66").Value = 23857.5
$ws.Range("L66").Value = 33000
$ws.Range("M66").Value = -20425.5
$ws.Range("N66").Value = -39864

$ws.Range("H122").Value = 2129.8215
$ws.Range("I122").Value = 2043.7693
$ws.Range("K122").Value = 6131.3079
$ws.Range("M122").Value = -3681.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 48684
$ws.Range("J53").Value = 48684
$ws.Range("L53").Value = 48684
$ws.Range("N53").Value = -49898

$ws.Range("H103").Value = 8066.3335
$ws.Range("J103").Value = 15999
$ws.Range("L103").Value = 15999
$ws.Range("N103").Value = -18343

$ws.Range("H107").Value = 685.6
$ws.Range("I107").Value = 762
$ws.Range("J107").Value = 660.13336
$ws.Range("K107").Value = 762
$ws.Range("L107").Value = 660.13336
$ws.Range("M107").Value = 1158
$ws.Range("N107").Value = -4500.13336

$ws.Range("H132").Value = 2212.182
$ws.Range("I132").Value = 1962.5128
$ws.Range("J132").Value = 4159.6
$ws.Range("K132").Value = 5887.538399999999
$ws.Range("L132").Value = 12478.8
$ws.Range("M132").Value = -3357.538399999999
$ws.Range("N132").Value = -17538.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 104008
$ws.Range("J88").Value = 104008
$ws.Range("L88").Value = 312024
$ws.Range("N88").Value = -312880

$ws.Range("H91").Value = 104008
$ws.Range("J91").Value = 104008
$ws.Range("L91").Value = 312024
$ws.Range("N91").Value = -314988

$ws.Range("H117").Value = 40110
$ws.Range("I117").Value = 330
$ws.Range("J117").Value = 43425
$ws.Range("K117").Value = 990
$ws.Range("L117").Value = 130275
$ws.Range("M117").Value = 2452
$ws.Range("N117").Value = -137159

$ws.Range("H132").Value = 1670.4445
$ws.Range("I132").Value = 925.38464
$ws.Range("J132").Value = 2091.5652
$ws.Range("K132").Value = 8328.46176
$ws.Range("L132").Value = 18824.0868
$ws.Range("M132").Value = -5798.46176
$ws.Range("N132").Value = -23884.0868

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 7777
$ws.Range("J117").Value = 7777
$ws.Range("L117").Value = 7777
$ws.Range("N117").Value = -14661

$ws.Range("H122").Value = 3241.5293
$ws.Range("I122").Value = 3444.963
$ws.Range("J122").Value = 2456.8572
$ws.Range("K122").Value = 10334.889
$ws.Range("L122").Value = 7370.571599999999
$ws.Range("M122").Value = -7884.889000000001
$ws.Range("N122").Value = -12270.5716

$ws.Range("H123").Value = 16887.105
$ws.Range("J123").Value = 16887.105
$ws.Range("L123").Value = 16887.105
$ws.Range("N123").Value = -21787.105

$ws.Range("H132").Value = 1943.5686
$ws.Range("I132").Value = 1632.7693
$ws.Range("J132").Value = 2953.6667
$ws.Range("K132").Value = 4898.3079
$ws.Range("L132").Value = 8861.000100000001
$ws.Range("M132").Value = -2368.3079
$ws.Range("N132").Value = -13921.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1614.2858
$ws.Range("I22").Value = 10000
$ws.Range("J22").Value = 216.66667
$ws.Range("K22").Value = 10000
$ws.Range("L22").Value = 216.66667
$ws.Range("M22").Value = -9705
$ws.Range("N22").Value = -806.6666700000001

$ws.Range("H27").Value = 1614.2858
$ws.Range("I27").Value = 10000
$ws.Range("J27").Value = 216.66667
$ws.Range("K27").Value = 10000
$ws.Range("L27").Value = 216.66667
$ws.Range("M27").Value = -9893
$ws.Range("N27").Value = -430.66667

$ws.Range("H136").Value = 3445.7021
$ws.Range("I136").Value = 3457.7954
$ws.Range("J136").Value = 3268.3333
$ws.Range("K136").Value = 10373.3862
$ws.Range("L136").Value = 9804.999899999999
$ws.Range("M136").Value = -7823.386200000001
$ws.Range("N136").Value = -14904.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 13980
$ws.Range("J30").Value = 13980
$ws.Range("L30").Value = 13980
$ws.Range("N30").Value = -14194

$ws.Range("H68").Value = 48723.25
$ws.Range("J68").Value = 48723.25
$ws.Range("L68").Value = 48723.25
$ws.Range("N68").Value = -50345.25

$ws.Range("H71").Value = 48723.25
$ws.Range("J71").Value = 48723.25
$ws.Range("L71").Value = 146169.75
$ws.Range("N71").Value = -154281.75

$ws.Range("H136").Value = 1431.9828
$ws.Range("I136").Value = 1427.0454
$ws.Range("J136").Value = 1447.5
$ws.Range("K136").Value = 4281.1362
$ws.Range("L136").Value = 4342.5
$ws.Range("M136").Value = -1731.1362
$ws.Range("N136").Value = -9442.5
